# Remove slide-level animation timing (the <p:timing> element) from every
# slide in the deck that currently has one. This mirrors an author action
# of clearing all animations on slides 2-7 (their MainSequence effects),
# which causes PowerPoint to drop the now-empty <p:timing> node on save.

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    $seq = $s.TimeLine.MainSequence
    for ($j = $seq.Count; $j -ge 1; $j--) {
        $seq.Item($j).Delete()
    }
}
